$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values that changed ---
$ws.Range("C5").Value = 434486
$ws.Range("C10").Value = 1694997

# --- Add new row 13 for year 2021 ---
# Copy formatting (style) from row 12's label cell (A12) onto A13 so the
# new label cell keeps the same bold/centered/bordered style (s="1").
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 1915236
$ws.Range("D13").Value = 24254118
$ws.Range("E13").Value = 6697618
